$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2053.4
$ws.Range("I18").Value = 983.6667
$ws.Range("J18").Value = 3658
$ws.Range("K18").Value = 983.6667
$ws.Range("L18").Value = 3658
$ws.Range("M18").Value = -699.6667
$ws.Range("N18").Value = -4226
$ws.Range("H40").Value = 3972
$ws.Range("J40").Value = 4444
$ws.Range("L40").Value = 4444
$ws.Range("N40").Value = -4794
$ws.Range("H43").Value = 1001
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 1002
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 1002
$ws.Range("M43").Value = -931
$ws.Range("N43").Value = -1140
$ws.Range("H51").Value = 10000
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").Value = $null
$ws.Range("H132").Value = 10525.4
$ws.Range("I132").Value = 10525.4
$ws.Range("K132").Value = 31576.2
$ws.Range("M132").Value = -29046.2
$ws.Range("H137").Value = 3397.4
$ws.Range("I137").Value = 2568.6428
$ws.Range("K137").Value = 7705.928400000001
$ws.Range("M137").Value = -5155.928400000001
$ws.Range("H138").Value = 2393.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4789.7856
$ws.Range("I32").Value = 4789.7856
$ws.Range("K32").Value = 4789.7856
$ws.Range("M32").Value = -4502.7856
$ws.Range("H45").Value = 4233
$ws.Range("I45").Value = 4249.5
$ws.Range("J45").Value = 4200
$ws.Range("K45").Value = 4249.5
$ws.Range("L45").Value = 4200
$ws.Range("M45").Value = -3872.5
$ws.Range("N45").Value = -4954
$ws.Range("H62").Value = 80000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 80000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = $null
$ws.Range("M62").Value = 80000
$ws.Range("N62").Value = -81248
$ws.Range("H65").Value = 80000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 80000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = $null
$ws.Range("M65").Value = 240000
$ws.Range("N65").Value = -246240
$ws.Range("H74").Value = 2415.25
$ws.Range("I74").Value = 1626.4615
$ws.Range("K74").Value = 1626.4615
$ws.Range("M74").Value = -752.4614999999999
$ws.Range("H77").Value = 2415.25
$ws.Range("I77").Value = 1626.4615
$ws.Range("K77").Value = 8132.307499999999
$ws.Range("M77").Value = -3764.307499999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 12499.5
$ws.Range("J18").Value = 12499.5
$ws.Range("L18").Value = 12499.5
$ws.Range("N18").Value = -13557.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 8060140
$ws.Range("I6").Value = 10075113
$ws.Range("K6").Value = 10075113
$ws.Range("M6").Value = -10075000
$ws.Range("H29").Value = 10021
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 10021
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = $null
$ws.Range("M29").Value = 10021
$ws.Range("N29").Value = -10607
$ws.Range("H41").Value = 1000
$ws.Range("I41").Value = 1000
$ws.Range("K41").Value = 1000
$ws.Range("M41").Value = -572
$ws.Range("H99").Value = 5599.25
$ws.Range("I99").Value = 5071.75
$ws.Range("J99").Value = 6126.75
$ws.Range("K99").Value = 5071.75
$ws.Range("L99").Value = 6126.75
$ws.Range("M99").Value = -3573.75
$ws.Range("N99").Value = -9122.75
$ws.Range("H126").Value = 5599.25
$ws.Range("I126").Value = 5071.75
$ws.Range("J126").Value = 6126.75
$ws.Range("K126").Value = 15215.25
$ws.Range("L126").Value = 18380.25
$ws.Range("M126").Value = -12745.25
$ws.Range("N126").Value = -23320.25
$ws.Range("H134").Value = 806.75
$ws.Range("I134").Value = 806.75
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 2420.25
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = $null
$ws.Range("N134").Value = 114.75

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1086.3334
$ws.Range("J34").Value = 1369.5
$ws.Range("L34").Value = 4108.5
$ws.Range("N34").Value = -4276.5
$ws.Range("J39").Value = 3000
$ws.Range("L39").Value = 9000
$ws.Range("N39").Value = -9588
$ws.Range("H48").Value = 524.5
$ws.Range("J48").Value = 944
$ws.Range("L48").Value = 2832
$ws.Range("N48").Value = -3332
$ws.Range("H55").Value = 2100
$ws.Range("J55").Value = 4000
$ws.Range("L55").Value = 12000
$ws.Range("N55").Value = -12354
$ws.Range("H92").Value = 1299.6
$ws.Range("J92").Value = 1324.5
$ws.Range("L92").Value = 3973.5
$ws.Range("N92").Value = -6469.5
$ws.Range("H113").Value = 603
$ws.Range("I113").Value = 603
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1809
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = $null
$ws.Range("N113").Value = 361
$ws.Range("H131").Value = 3138.5715
$ws.Range("I131").Value = 1161.6666
$ws.Range("J131").Value = 15000
$ws.Range("K131").Value = 3484.9998
$ws.Range("L131").Value = 45000
$ws.Range("M131").Value = 1555.0002
$ws.Range("N131").Value = -55080

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 4995
$ws.Range("I5").Value = 4995
$ws.Range("K5").Value = 4995
$ws.Range("M5").Value = -4883
$ws.Range("H18").Value = 2000740
$ws.Range("I18").Value = 10000000
$ws.Range("K18").Value = 10000000
$ws.Range("M18").Value = -9999707
$ws.Range("H113").Value = 490
$ws.Range("I113").Value = 490
$ws.Range("K113").Value = 490
$ws.Range("M113").Value = 1680

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5031.6665
$ws.Range("I7").Value = 5031.6665
$ws.Range("K7").Value = 5031.6665
$ws.Range("M7").Value = -4919.6665
$ws.Range("H63").Value = 40000
$ws.Range("J63").Value = 40000
$ws.Range("L63").Value = 40000
$ws.Range("N63").Value = -41498
$ws.Range("H66").Value = 40000
$ws.Range("J66").Value = 40000
$ws.Range("L66").Value = 120000
$ws.Range("N66").Value = -127488
$ws.Range("H122").Value = 4331.25
$ws.Range("J122").Value = 3250
$ws.Range("L122").Value = 9750
$ws.Range("N122").Value = -14650
$ws.Range("H126").Value = 5031.6665
$ws.Range("I126").Value = 5031.6665
$ws.Range("K126").Value = 15094.9995
$ws.Range("M126").Value = -12624.9995
$ws.Range("H132").Value = 1288
$ws.Range("I132").Value = 1288
$ws.Range("K132").Value = 3864
$ws.Range("M132").Value = -1334
$ws.Range("H136").Value = 5688.75
$ws.Range("I136").Value = 5651.6665
$ws.Range("J136").Value = 5800
$ws.Range("K136").Value = 16954.9995
$ws.Range("L136").Value = 17400
$ws.Range("M136").Value = -14404.9995
$ws.Range("N136").Value = -22500

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 668
$ws.Range("I14").Value = 668
$ws.Range("K14").Value = 668
$ws.Range("M14").Value = -500
$ws.Range("H132").Value = 1343.091
$ws.Range("I132").Value = 1343.091
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4029.273
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = -1499.273
